# Update "想去人数" (want-to-go count) values in column F on the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets to reflect the
# latest scraped numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 74
$ws1.Range("F7").Value = 592
$ws1.Range("F8").Value = 111
$ws1.Range("F9").Value = 8736
$ws1.Range("F10").Value = 808
$ws1.Range("F13").Value = 985
$ws1.Range("F14").Value = 111
$ws1.Range("F16").Value = 6
$ws1.Range("F18").Value = 261
$ws1.Range("F21").Value = 1031

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 74
$ws4.Range("F9").Value = 592
$ws4.Range("F10").Value = 111
$ws4.Range("F11").Value = 8736
$ws4.Range("F12").Value = 808
$ws4.Range("F15").Value = 985
$ws4.Range("F16").Value = 111
$ws4.Range("F18").Value = 6
$ws4.Range("F20").Value = 261
$ws4.Range("F23").Value = 1031
